# Update weekly Haba (bean) price records: each row's date / quality /
# volume / price / origin fields are reshuffled among the existing rows,
# per the "Fruta / hortaliza, semanal" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44466
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("O2").Value = 'Región de O''Higgins'
$ws.Range("P2").Value = 460
$ws.Range("D3").Value = 44165
$ws.Range("J3").Value = 38
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8500
$ws.Range("M3").Value = 8263
$ws.Range("O3").Value = 'Región del Maule'
$ws.Range("P3").Value = 331
$ws.Range("D4").Value = 44488
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8500
$ws.Range("O4").Value = 'Región del Maule'
$ws.Range("P4").Value = 340
$ws.Range("D5").Value = 44518
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6500
$ws.Range("O5").Value = 'Provincia de Diguillín'
$ws.Range("P5").Value = 260
$ws.Range("D6").Value = 44519
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 6000
$ws.Range("M6").Value = 6500
$ws.Range("P6").Value = 260
$ws.Range("D7").Value = 44516
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 300
$ws.Range("D8").Value = 44515
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 7500
$ws.Range("P8").Value = 300
$ws.Range("D9").Value = 44536
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6750
$ws.Range("O9").Value = 'Provincia de Diguillín'
$ws.Range("P9").Value = 270
$ws.Range("D10").Value = 44167
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8500
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 340
$ws.Range("D11").Value = 44484
$ws.Range("J11").Value = 30
$ws.Range("K11").Value = 8500
$ws.Range("M11").Value = 8750
$ws.Range("P11").Value = 350
$ws.Range("D12").Value = 44537
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6750
$ws.Range("O12").Value = 'Provincia de Diguillín'
$ws.Range("P12").Value = 270
$ws.Range("D13").Value = 44522
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 6000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6500
$ws.Range("O13").Value = 'Provincia de Diguillín'
$ws.Range("P13").Value = 260
$ws.Range("D14").Value = 44476
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 7500
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7750
$ws.Range("P14").Value = 310
$ws.Range("D15").Value = 44495
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8500
$ws.Range("O15").Value = 'Región del Maule'
$ws.Range("P15").Value = 340
$ws.Range("D16").Value = 44159
$ws.Range("J16").Value = 42
$ws.Range("K16").Value = 6500
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 6738
$ws.Range("P16").Value = 270
$ws.Range("D17").Value = 44509
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = 8500
$ws.Range("O17").Value = 'Región del Maule'
$ws.Range("P17").Value = 340
$ws.Range("D18").Value = 44523
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6500
$ws.Range("O18").Value = 'Provincia de Diguillín'
$ws.Range("P18").Value = 260
$ws.Range("D19").Value = 44566
$ws.Range("K19").Value = 7000
$ws.Range("L19").Value = 7500
$ws.Range("M19").Value = 7250
$ws.Range("P19").Value = 290
$ws.Range("D20").Value = 44161
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = 6500
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 6764
$ws.Range("O20").Value = 'Región de O''Higgins'
$ws.Range("P20").Value = 271
$ws.Range("D21").Value = 44526
$ws.Range("K21").Value = 6000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 6500
$ws.Range("P21").Value = 260
$ws.Range("D22").Value = 44166
$ws.Range("J22").Value = 56
$ws.Range("K22").Value = 7500
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 7804
$ws.Range("O22").Value = 'Región de O''Higgins'
$ws.Range("P22").Value = 312
$ws.Range("D23").Value = 44533
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 6500
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6750
$ws.Range("O23").Value = 'Provincia de Diguillín'
$ws.Range("P23").Value = 270
$ws.Range("D24").Value = 44517
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 6000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 6500
$ws.Range("O24").Value = 'Provincia de Diguillín'
$ws.Range("P24").Value = 260
$ws.Range("D25").Value = 44524
$ws.Range("J25").Value = 100
$ws.Range("D26").Value = 44530
$ws.Range("J26").Value = 100
$ws.Range("D27").Value = 44160
$ws.Range("K27").Value = 6500
$ws.Range("M27").Value = 6688
$ws.Range("O27").Value = 'Región de O''Higgins'
$ws.Range("P27").Value = 268
$ws.Range("D28").Value = 44162
$ws.Range("J28").Value = 80
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 7562
$ws.Range("P28").Value = 302
$ws.Range("D29").Value = 44473
$ws.Range("K29").Value = 9500
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 9750
$ws.Range("P29").Value = 390
$ws.Range("D30").Value = 44512
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 8000
$ws.Range("M30").Value = 7500
$ws.Range("P30").Value = 300
$ws.Range("D31").Value = 44540
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 6500
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = 6750
$ws.Range("O31").Value = 'Provincia de Diguillín'
$ws.Range("P31").Value = 270
$ws.Range("D32").Value = 44482
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 8000
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = 8500
$ws.Range("O32").Value = 'Región del Maule'
$ws.Range("P32").Value = 340
$ws.Range("D33").Value = 44553
$ws.Range("K33").Value = 6500
$ws.Range("L33").Value = 7000
$ws.Range("M33").Value = 6750
$ws.Range("P33").Value = 270
$ws.Range("D34").Value = 44511
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = 7500
$ws.Range("P34").Value = 300
$ws.Range("D35").Value = 44487
$ws.Range("J35").Value = 30
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = 8000
$ws.Range("P35").Value = 320
$ws.Range("D36").Value = 44487
$ws.Range("I36").Value = 'Segunda'
$ws.Range("J36").Value = 30
$ws.Range("K36").Value = 9000
$ws.Range("L36").Value = 9000
$ws.Range("M36").Value = 9000
$ws.Range("O36").Value = 'Región del Maule'
$ws.Range("P36").Value = 360
$ws.Range("D37").Value = 44529
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 6000
$ws.Range("L37").Value = 7000
$ws.Range("M37").Value = 6500
$ws.Range("O37").Value = 'Provincia de Diguillín'
$ws.Range("P37").Value = 260
$ws.Range("D38").Value = 44504
$ws.Range("J38").Value = 60
$ws.Range("K38").Value = 8000
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = 8500
$ws.Range("O38").Value = 'Región del Maule'
$ws.Range("P38").Value = 340
$ws.Range("D39").Value = 44491
$ws.Range("D40").Value = 44489
$ws.Range("J40").Value = 60
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = 8500
$ws.Range("O40").Value = 'Región del Maule'
$ws.Range("P40").Value = 340
